# Fruta / hortaliza, semanal
# Re-associates the (Fecha, Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) tuple for each data row (2-44)
# with the tuple that belongs to a different row, per the weekly re-shuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (the row whose D/L/M/N/O/P/S tuple
# should end up in the target row). Row 16 maps to itself (no change).
$mapping = @{
    2=42; 3=28; 4=19; 5=41; 6=4; 7=11; 8=25; 9=10; 10=6;
    11=13; 12=24; 13=35; 14=26; 15=17; 16=16; 17=39; 18=38; 19=44;
    20=36; 21=14; 22=3; 23=9; 24=22; 25=29; 26=37; 27=31;
    28=30; 29=20; 30=12; 31=5; 32=34; 33=32; 34=15; 35=23;
    36=43; 37=21; 38=2; 39=8; 40=27; 41=33; 42=40; 43=18; 44=7
}

$firstRow = 2
$lastRow = 44

# Snapshot the current (before) values for the columns that move, keyed by row.
$snapD = @{}
$snapL = @{}
$snapM = @{}
$snapN = @{}
$snapO = @{}
$snapP = @{}
$snapS = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, 4).Value2   # D: Fecha
    $snapL[$r] = $ws.Cells.Item($r, 12).Value2  # L: Calidad
    $snapM[$r] = $ws.Cells.Item($r, 13).Value2  # M: Volumen
    $snapN[$r] = $ws.Cells.Item($r, 14).Value2  # N: Precio minimo
    $snapO[$r] = $ws.Cells.Item($r, 15).Value2  # O: Precio maximo
    $snapP[$r] = $ws.Cells.Item($r, 16).Value2  # P: Precio promedio ponderado
    $snapS[$r] = $ws.Cells.Item($r, 19).Value2  # S: Precio $/Kg
}

# Now write the new values based on the mapping, using the pre-captured snapshot
# so that reassignments never read already-overwritten data.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $mapping[$r]
    if ($src -eq $r) { continue }

    $ws.Cells.Item($r, 4).Value = $snapD[$src]
    $ws.Cells.Item($r, 12).Value = $snapL[$src]
    $ws.Cells.Item($r, 13).Value = $snapM[$src]
    $ws.Cells.Item($r, 14).Value = $snapN[$src]
    $ws.Cells.Item($r, 15).Value = $snapO[$src]
    $ws.Cells.Item($r, 16).Value = $snapP[$src]
    $ws.Cells.Item($r, 19).Value = $snapS[$src]
}
